$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 2).Value = 26.67882461893246
$ws.Cells.Item(2, 4).Value = 9.339034718821546
$ws.Cells.Item(2, 5).Value = 10.81369437765513
$ws.Cells.Item(2, 6).Value = 62.72601528416051
$ws.Cells.Item(2, 7).Value = 3.805974943425642
$ws.Cells.Item(2, 10).Value = 8.32751604154571
$ws.Cells.Item(2, 12).Value = 18.85716188208018
$ws.Cells.Item(2, 14).Value = 20.96866967226643
$ws.Cells.Item(3, 2).Value = 26.43006989693861
$ws.Cells.Item(3, 4).Value = 9.011715909547659
$ws.Cells.Item(3, 5).Value = 10.1469518066817
$ws.Cells.Item(3, 6).Value = 62.41385843880992
$ws.Cells.Item(3, 7).Value = 3.813780529659272
$ws.Cells.Item(3, 10).Value = 8.157976105173558
$ws.Cells.Item(3, 12).Value = 18.66440368713422
$ws.Cells.Item(3, 14).Value = 21.06179260090708
$ws.Cells.Item(4, 2).Value = 26.28393072089948
$ws.Cells.Item(4, 4).Value = 8.808833168457578
$ws.Cells.Item(4, 5).Value = 9.731616932367251
$ws.Cells.Item(4, 6).Value = 62.24310805665611
$ws.Cells.Item(4, 7).Value = 3.818806386003851
$ws.Cells.Item(4, 10).Value = 8.051036646518803
$ws.Cells.Item(4, 12).Value = 18.55079640002241
$ws.Cells.Item(4, 14).Value = 21.12116226465107
$ws.Cells.Item(5, 2).Value = 26.22608913408453
$ws.Cells.Item(5, 4).Value = 8.725819778496342
$ws.Cells.Item(5, 5).Value = 9.589939928495159
$ws.Cells.Item(5, 6).Value = 62.17879871260058
$ws.Cells.Item(5, 7).Value = 3.820913441112234
$ws.Cells.Item(5, 10).Value = 8.006768386060555
$ws.Cells.Item(5, 12).Value = 18.50573504567653
$ws.Cells.Item(5, 14).Value = 21.14591030606808
$ws.Cells.Item(6, 2).Value = 26.21658940159335
$ws.Cells.Item(6, 4).Value = 8.712019480010055
$ws.Cells.Item(6, 5).Value = 9.566266619117556
$ws.Cells.Item(6, 6).Value = 62.16843888463188
$ws.Cells.Item(6, 7).Value = 3.82126688802735
$ws.Cells.Item(6, 10).Value = 7.999376780320429
$ws.Cells.Item(6, 12).Value = 18.4983283451065
$ws.Cells.Item(6, 14).Value = 21.15005328064028
$ws.Cells.Item(7, 2).Value = 26.28314365340444
$ws.Cells.Item(7, 4).Value = 8.807714789446374
$ws.Cells.Item(7, 5).Value = 9.729716289842777
$ws.Cells.Item(7, 6).Value = 62.24221939735099
$ws.Cells.Item(7, 7).Value = 3.818834563237241
$ws.Cells.Item(7, 10).Value = 8.050442387460691
$ws.Cells.Item(7, 12).Value = 18.55018363626053
$ws.Cells.Item(7, 14).Value = 21.12149377622137
$ws.Cells.Item(8, 2).Value = 26.59171723479226
$ws.Cells.Item(8, 4).Value = 9.226653509386882
$ws.Cells.Item(8, 5).Value = 10.58865028963367
$ws.Cells.Item(8, 6).Value = 62.61405337404843
$ws.Cells.Item(8, 7).Value = 3.80861811514685
$ws.Cells.Item(8, 10).Value = 8.269662882051973
$ws.Cells.Item(8, 12).Value = 18.78973893460699
$ws.Cells.Item(8, 14).Value = 21.00032553572586
$ws.Cells.Item(9, 2).Value = 27.24671316856629
$ws.Cells.Item(9, 4).Value = 10.02733236418452
$ws.Cells.Item(9, 5).Value = 12.12274626110275
$ws.Cells.Item(9, 6).Value = 63.5081849059318
$ws.Cells.Item(9, 7).Value = 3.790418456540242
$ws.Cells.Item(9, 10).Value = 8.675990751000469
$ws.Cells.Item(9, 12).Value = 19.29525599399884
$ws.Cells.Item(9, 14).Value = 20.77995346739295
$ws.Cells.Item(10, 2).Value = 27.75486008365623
$ws.Cells.Item(10, 4).Value = 10.59601620232753
$ws.Cells.Item(10, 5).Value = 13.13689848105807
$ws.Cells.Item(10, 6).Value = 64.26380895169179
$ws.Cells.Item(10, 7).Value = 3.778143841300367
$ws.Cells.Item(10, 10).Value = 8.958900769339229
$ws.Cells.Item(10, 12).Value = 19.68574797920347
$ws.Cells.Item(10, 14).Value = 20.62833688780823
$ws.Cells.Item(11, 2).Value = 27.99108857657904
$ws.Cells.Item(11, 4).Value = 10.84925485393607
$ws.Cells.Item(11, 5).Value = 13.57384753372484
$ws.Cells.Item(11, 6).Value = 64.62843359293969
$ws.Cells.Item(11, 7).Value = 3.77279318462552
$ws.Cells.Item(11, 10).Value = 9.083977006850514
$ws.Cells.Item(11, 12).Value = 19.86693168091395
$ws.Cells.Item(11, 14).Value = 20.56154941862882
$ws.Cells.Item(12, 2).Value = 28.08120311697694
$ws.Cells.Item(12, 4).Value = 10.94427813529611
$ws.Cells.Item(12, 5).Value = 13.73581326783227
$ws.Cells.Item(12, 6).Value = 64.76945205565291
$ws.Cells.Item(12, 7).Value = 3.770800175636073
$ws.Cells.Item(12, 10).Value = 9.130801579346032
$ws.Cells.Item(12, 12).Value = 19.93599910554344
$ws.Cells.Item(12, 14).Value = 20.53656902484551
$ws.Cells.Item(13, 2).Value = 28.06176709428535
$ws.Cells.Item(13, 4).Value = 10.92385327077542
$ws.Cells.Item(13, 5).Value = 13.70108640264661
$ws.Cells.Item(13, 6).Value = 64.73895134832723
$ws.Cells.Item(13, 7).Value = 3.771227936289023
$ws.Cells.Item(13, 10).Value = 9.120741322488104
$ws.Cells.Item(13, 12).Value = 19.92110471635336
$ws.Cells.Item(13, 14).Value = 20.54193524200424
$ws.Cells.Item(14, 2).Value = 27.99848952603993
$ws.Cells.Item(14, 4).Value = 10.85709043193612
$ws.Cells.Item(14, 5).Value = 13.5872425759343
$ws.Cells.Item(14, 6).Value = 64.63997665436847
$ws.Cells.Item(14, 7).Value = 3.7726285558755
$ws.Cells.Item(14, 10).Value = 9.087840196587072
$ws.Cells.Item(14, 12).Value = 19.87260504297178
$ws.Cells.Item(14, 14).Value = 20.55948806213614
$ws.Cells.Item(15, 2).Value = 27.95981406952948
$ws.Cells.Item(15, 4).Value = 10.81608018946176
$ws.Cells.Item(15, 5).Value = 13.51705483610121
$ws.Cells.Item(15, 6).Value = 64.57973315557265
$ws.Cells.Item(15, 7).Value = 3.77349078474767
$ws.Cells.Item(15, 10).Value = 9.067616582372548
$ws.Cells.Item(15, 12).Value = 19.84295554294293
$ws.Cells.Item(15, 14).Value = 20.57028001414993
$ws.Cells.Item(16, 2).Value = 27.73951812429666
$ws.Cells.Item(16, 4).Value = 10.57934887592401
$ws.Cells.Item(16, 5).Value = 13.10785279466644
$ws.Cells.Item(16, 6).Value = 64.24039515734636
$ws.Cells.Item(16, 7).Value = 3.778498180159712
$ws.Cells.Item(16, 10).Value = 8.950652272974672
$ws.Cells.Item(16, 12).Value = 19.67397405212906
$ws.Cells.Item(16, 14).Value = 20.63274525398418
$ws.Cells.Item(17, 2).Value = 27.60562329827285
$ws.Cells.Item(17, 4).Value = 10.43265798456819
$ws.Cells.Item(17, 5).Value = 12.85057927697322
$ws.Cells.Item(17, 6).Value = 64.03753110768179
$ws.Cells.Item(17, 7).Value = 3.781629513144675
$ws.Cells.Item(17, 10).Value = 8.877957473028845
$ws.Cells.Item(17, 12).Value = 19.57118027600548
$ws.Cells.Item(17, 14).Value = 20.67162251132657
$ws.Cells.Item(18, 2).Value = 27.52909193994433
$ws.Cells.Item(18, 4).Value = 10.34777603717381
$ws.Cells.Item(18, 5).Value = 12.70030547092789
$ws.Cells.Item(18, 6).Value = 63.92281917104638
$ws.Cells.Item(18, 7).Value = 3.783452537538243
$ws.Cells.Item(18, 10).Value = 8.835805528742014
$ws.Cells.Item(18, 12).Value = 19.51239307415565
$ws.Cells.Item(18, 14).Value = 20.69418942700968
$ws.Cells.Item(19, 2).Value = 27.50326451186872
$ws.Cells.Item(19, 4).Value = 10.31895185776544
$ws.Cells.Item(19, 5).Value = 12.64903057809555
$ws.Cells.Item(19, 6).Value = 63.88431978206675
$ws.Cells.Item(19, 7).Value = 3.78407356544274
$ws.Cells.Item(19, 10).Value = 8.821475826249223
$ws.Cells.Item(19, 12).Value = 19.49254823574639
$ws.Cells.Item(19, 14).Value = 20.70186563794232
$ws.Cells.Item(20, 2).Value = 27.61982730601824
$ws.Cells.Item(20, 4).Value = 10.44832685292966
$ws.Cells.Item(20, 5).Value = 12.87820413375111
$ws.Cells.Item(20, 6).Value = 64.05892285210292
$ws.Cells.Item(20, 7).Value = 3.78129390672314
$ws.Cells.Item(20, 10).Value = 8.885731266539215
$ws.Cells.Item(20, 12).Value = 19.58208834725216
$ws.Cells.Item(20, 14).Value = 20.66746269183955
$ws.Cells.Item(21, 2).Value = 28.01705832501455
$ws.Cells.Item(21, 4).Value = 10.87672462818398
$ws.Cells.Item(21, 5).Value = 13.62077608958651
$ws.Cells.Item(21, 6).Value = 64.66896856328864
$ws.Cells.Item(21, 7).Value = 3.772216262609769
$ws.Cells.Item(21, 10).Value = 9.097518824468111
$ws.Cells.Item(21, 12).Value = 19.88683859065893
$ws.Cells.Item(21, 14).Value = 20.55432397212291
$ws.Cells.Item(22, 2).Value = 28.28048599202808
$ws.Cells.Item(22, 4).Value = 11.15158717344568
$ws.Cells.Item(22, 5).Value = 14.08571092005673
$ws.Cells.Item(22, 6).Value = 65.08479034446918
$ws.Cells.Item(22, 7).Value = 3.766476662024511
$ws.Cells.Item(22, 10).Value = 9.232784579693423
$ws.Cells.Item(22, 12).Value = 20.08865010889969
$ws.Cells.Item(22, 14).Value = 20.48218990408585
$ws.Cells.Item(23, 2).Value = 28.13956365318332
$ws.Cells.Item(23, 4).Value = 11.00538265773428
$ws.Cells.Item(23, 5).Value = 13.83942701862163
$ws.Cells.Item(23, 6).Value = 64.86131346540776
$ws.Cells.Item(23, 7).Value = 3.7695224384621
$ws.Cells.Item(23, 10).Value = 9.16088457900438
$ws.Cells.Item(23, 12).Value = 19.98071551658371
$ws.Cells.Item(23, 14).Value = 20.5205248698728
$ws.Cells.Item(24, 2).Value = 27.61340427794618
$ws.Cells.Item(24, 4).Value = 10.44124465579708
$ws.Cells.Item(24, 5).Value = 12.86572229957634
$ws.Cells.Item(24, 6).Value = 64.04924566794786
$ws.Cells.Item(24, 7).Value = 3.781445563406949
$ws.Cells.Item(24, 10).Value = 8.882217850461229
$ws.Cells.Item(24, 12).Value = 19.57715583768384
$ws.Cells.Item(24, 14).Value = 20.66934267360653
$ws.Cells.Item(25, 2).Value = 27.06452980080271
$ws.Cells.Item(25, 4).Value = 9.813713248835775
$ws.Cells.Item(25, 5).Value = 11.72759188349091
$ws.Cells.Item(25, 6).Value = 63.24879787576722
$ws.Cells.Item(25, 7).Value = 3.795147770689703
$ws.Cells.Item(25, 10).Value = 8.568726087353804
$ws.Cells.Item(25, 12).Value = 19.15494558358967
$ws.Cells.Item(25, 14).Value = 20.83774690282723
